# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before the "总计" (Total) sheet,
#    and fill it with the per-fund holdings detail for the new quarter.
# 2. Update the "总计" (Total) summary sheet by inserting a new top data row
#    for "2022-Q1" and shifting the previously existing rows down.

$wb = $excel.ActiveWorkbook

$templateSheet = $wb.Worksheets.Item("2021-Q3")

# ------------------------------------------------------------------
# 1. Create the new "2022-Q1" detail sheet, positioned just before "总计"
# ------------------------------------------------------------------
$q1 = $wb.Worksheets.Add($wb.Worksheets.Item("总计"))
$q1.Name = "2022-Q1"

# NOTE: after the insertion above, re-resolve the "总计" sheet by name
# instead of reusing any handle obtained before the insert -- worksheet
# handles in this automation layer follow sheet *position*, and "总计"
# moved from position 4 to position 5 when "2022-Q1" was inserted before it.
$totalSheet = $wb.Worksheets.Item("总计")

# Copy header + index-column formatting from an existing detail sheet so the
# new sheet matches the look (borders/bold/centered) of the others.
$templateSheet.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats

$templateSheet.Range("A2").Copy()
$q1.Range("A2:A6").PasteSpecial(-4122)   # xlPasteFormats

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Columns B..G hold text-like values (fund codes with leading zeros, and
# decimal figures that must stay as text rather than be coerced to numbers),
# so format that block as Text before writing the values.
$q1.Range("B2:G6").NumberFormat = "@"

# Row 2
$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "000988"
$q1.Range("C2").Value = "嘉实全球互联网股票 - 人民币QDII"
$q1.Range("D2").Value = "13.21"
$q1.Range("E2").Value = "85.88"
$q1.Range("F2").Value = "3.83"
$q1.Range("G2").Value = "0.5059"
$q1.Range("H2").Value = 9

# Row 3
$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "000989"
$q1.Range("C3").Value = "嘉实全球互联网股票 - 美元现汇QDII"
$q1.Range("D3").Value = "13.21"
$q1.Range("E3").Value = "85.88"
$q1.Range("F3").Value = "3.83"
$q1.Range("G3").Value = "0.5059"
$q1.Range("H3").Value = 9

# Row 4
$q1.Range("A4").Value = 2
$q1.Range("B4").Value = "000990"
$q1.Range("C4").Value = "嘉实全球互联网股票 - 美元现钞QDII"
$q1.Range("D4").Value = "13.21"
$q1.Range("E4").Value = "85.88"
$q1.Range("F4").Value = "3.83"
$q1.Range("G4").Value = "0.5059"
$q1.Range("H4").Value = 9

# Row 5
$q1.Range("A5").Value = 3
$q1.Range("B5").Value = "003721"
$q1.Range("C5").Value = "易方达标普信息科技指数（QDII-LOF）美元"
$q1.Range("D5").Value = "6.31"
$q1.Range("E5").Value = "93.58"
$q1.Range("F5").Value = "1.88"
$q1.Range("G5").Value = "0.1186"
$q1.Range("H5").Value = 8

# Row 6
$q1.Range("A6").Value = 4
$q1.Range("B6").Value = "161128"
$q1.Range("C6").Value = "易方达标普信息科技指数（QDII-LOF）人民币"
$q1.Range("D6").Value = "6.31"
$q1.Range("E6").Value = "93.58"
$q1.Range("F6").Value = "1.88"
$q1.Range("G6").Value = "0.1186"
$q1.Range("H6").Value = 8

# ------------------------------------------------------------------
# 2. Update the "总计" sheet: push existing rows down one slot and add the
#    new "2022-Q1" row at the top of the data (row 2).
# ------------------------------------------------------------------
$b4 = $totalSheet.Range("B4").Value2
$c4 = $totalSheet.Range("C4").Value2
$d4 = $totalSheet.Range("D4").Value2

$b3 = $totalSheet.Range("B3").Value2
$c3 = $totalSheet.Range("C3").Value2
$d3 = $totalSheet.Range("D3").Value2

$b2 = $totalSheet.Range("B2").Value2
$c2 = $totalSheet.Range("C2").Value2
$d2 = $totalSheet.Range("D2").Value2

$totalSheet.Range("B5").Value = $b4
$totalSheet.Range("C5").Value = $c4
$totalSheet.Range("D5").Value = $d4

$totalSheet.Range("B4").Value = $b3
$totalSheet.Range("C4").Value = $c3
$totalSheet.Range("D4").Value = $d3

$totalSheet.Range("B3").Value = $b2
$totalSheet.Range("C3").Value = $c2
$totalSheet.Range("D3").Value = $d2

$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 5
$totalSheet.Range("D2").Value = 1.75

# Re-number / re-style the index column A2:A5
$totalSheet.Range("A2").Copy()
$totalSheet.Range("A2:A5").PasteSpecial(-4122)   # xlPasteFormats

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
